$wb = $excel.ActiveWorkbook

# --- Backlog sheet: add owner/status/points details for a few user stories ---
$backlog = $wb.Worksheets.Item("Backlog")

# Row 31 (US30 - List living married): status + points
$backlog.Range("E31").Value = "in work"
$backlog.Range("F31").Value = 3

# Row 32 (US31 - List living single): status + points
$backlog.Range("E32").Value = "in work"
$backlog.Range("F32").Value = 3

# Row 36 (US35 - List recent births): sprint, owner, status, points
$backlog.Range("A36").Value = 3
$backlog.Range("D36").Value = "jj"
$backlog.Range("E36").Value = "in work"
$backlog.Range("F36").Value = 3

# Row 37 (US36 - List recent deaths): sprint, owner, status, points
$backlog.Range("A37").Value = 3
$backlog.Range("D37").Value = "jj"
$backlog.Range("E37").Value = "in work"
$backlog.Range("F37").Value = 3

# --- Stories sheet: row 28 description needed a taller row to fit wrapped text ---
$stories = $wb.Worksheets.Item("Stories")
$stories.Rows.Item(28).RowHeight = 31.5


# --- View state: user ended up with Sprint2 active, Backlog scrolled further down ---
$backlog.Range("G37").Select()
$backlogWindow = $excel.ActiveWindow
$backlogWindow.ScrollRow = 19

$sprint2 = $wb.Worksheets.Item("Sprint2")
$sprint2.Activate()
